$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.504.20'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '3.252.72'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.12'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.22'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.595'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.66'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').Value = '3.816.44'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.22'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '68.470.60'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '3.240.07'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.48'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '394.11'
$ws.Range('E20').Value = '  +4.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.66'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.64'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000119'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('E26').Value = '  +4.33%  '
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.98'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.11'
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.39'
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('E37').Value = '  +4.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.822'
$ws.Range('E38').Value = '  -3.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.59'
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.59'
$ws.Range('E41').Value = '  -3.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '41.44'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0689'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.47'
$ws.Range('E44').Value = '  -6.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '342.63'
$ws.Range('E45').Value = '  -3.26%  '
$ws.Range('D46').Value = '2.586.36'
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '31.82'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.102'
$ws.Range('E51').Value = '  -1.23%  '

# Restore default (General) styling on the cells we forced to
# text above, so we do not leave a stray number-format on them.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
